# Swap the deck's theme colours: the slide master currently carries the
# "Integral" colour scheme (persisted to ppt/theme/theme2.xml); the edit
# repoints it at the stock "Office Theme" colour scheme (the palette that
# used to live untouched in ppt/theme/theme1.xml, used only by the Notes
# Master). Font scheme / format scheme are already identical between the
# two themes, so only the 12 theme colour slots need to move.

function RGBVal([int]$r, [int]$g, [int]$b) {
    return $r + ($g * 256) + ($b * 65536)
}

# Target palette = the built-in "Office Theme" colour scheme.
$officeTheme = @{
    dk1      = RGBVal 0x00 0x00 0x00
    lt1      = RGBVal 0xFF 0xFF 0xFF
    dk2      = RGBVal 0x44 0x54 0x6A
    lt2      = RGBVal 0xE7 0xE6 0xE6
    accent1  = RGBVal 0x5B 0x9B 0xD5
    accent2  = RGBVal 0xED 0x7D 0x31
    accent3  = RGBVal 0xA5 0xA5 0xA5
    accent4  = RGBVal 0xFF 0xC0 0x00
    accent5  = RGBVal 0x44 0x72 0xC4
    accent6  = RGBVal 0x70 0xAD 0x47
    hlink    = RGBVal 0x05 0x63 0xC1
    folHlink = RGBVal 0x95 0x4F 0x72
}

# ThemeColorScheme.Item() order matches the OOXML <a:clrScheme> child order.
$order = @("dk1","lt1","dk2","lt2","accent1","accent2","accent3","accent4","accent5","accent6","hlink","folHlink")

$p = $ppt.ActivePresentation
$master = $p.Slides.Item(1).Master
$colorScheme = $master.Theme.ThemeColorScheme

for ($i = 1; $i -le $order.Count; $i++) {
    $slotName = $order[$i - 1]
    $colorScheme.Item($i).RGB = $officeTheme[$slotName]
}
